$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price records arrived for Chirimoya (Vega Modelo de Temuco).
# Insert two new rows right above existing row 173; this pushes the old
# rows 173-201 down to 175-203 (matching the dimension growing to T203).
$ws.Rows("173:174").Insert()

# --- New row 173 ---
$ws.Range("A173").Value = 10
$ws.Range("B173").Value = 'Vega Modelo de Temuco'
$ws.Range("C173").Value = 'La Araucanía'
$ws.Range("D173").Value = 45173
$ws.Range("E173").Value = 9
$ws.Range("F173").Value = 'Fruta'
$ws.Range("G173").Value = 100107
$ws.Range("H173").Value = 'Otros'
$ws.Range("I173").Value = 100107002
$ws.Range("J173").Value = 'Chirimoya'
$ws.Range("K173").Value = 'Cultivar IV Región'
$ws.Range("L173").Value = 'Especial'
$ws.Range("M173").Value = 80
$ws.Range("N173").Value = 3500
$ws.Range("O173").Value = 3500
$ws.Range("P173").Value = 3500
$ws.Range("Q173").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R173").Value = 'Provincia del Elquí'
$ws.Range("S173").Value = 3500
$ws.Range("T173").Value = 1

# --- New row 174 ---
$ws.Range("A174").Value = 10
$ws.Range("B174").Value = 'Vega Modelo de Temuco'
$ws.Range("C174").Value = 'La Araucanía'
$ws.Range("D174").Value = 45173
$ws.Range("E174").Value = 9
$ws.Range("F174").Value = 'Fruta'
$ws.Range("G174").Value = 100107
$ws.Range("H174").Value = 'Otros'
$ws.Range("I174").Value = 100107002
$ws.Range("J174").Value = 'Chirimoya'
$ws.Range("K174").Value = 'Cultivar IV Región'
$ws.Range("L174").Value = 'Primera'
$ws.Range("M174").Value = 100
$ws.Range("N174").Value = 3000
$ws.Range("O174").Value = 3000
$ws.Range("P174").Value = 3000
$ws.Range("Q174").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R174").Value = 'Provincia del Elquí'
$ws.Range("S174").Value = 3000
$ws.Range("T174").Value = 1
